$p = $ppt.ActivePresentation

# --- Slide 2 ("3 things I learned in core") ---
$s2 = $p.Slides.Item(2)

# Title 1: "T" -> "3 things I learned in core"
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "3 things I learned in core"

# --- Slide 3 ("3 Things I leaned in Tech") ---
$s3 = $p.Slides.Item(3)

# Title 1: (empty) -> "3 Things I leaned in Tech"
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "3 Things I leaned in Tech"

# Content Placeholder 2: (empty) -> 4 paragraphs, second one indented one level
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "How to make a close button which requires a confirmation period"
[void]$tr3.InsertAfter("`rMy own invention")
[void]$tr3.InsertAfter("`rHow to make multiple buttons be bonded to a single event")
[void]$tr3.InsertAfter("`rHow to integrate ")
$tr3.Paragraphs(2).IndentLevel = 2
